# Adds a new weekly price record for "Haba" (Femacal de La Calera) on top of
# the existing list. The new record is inserted as row 23 (the list is kept
# in reverse-chronological / insertion order), pushing the previous rows
# 23-83 down to 24-84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 23; this shifts rows 23:83 down to
# 24:84 (and carries the date-format style already used by column D along
# with it), growing the used range from A1:R83 to A1:R84.
$ws.Rows(23).Insert()

# Populate the newly inserted row 23 with the new observation.
$ws.Cells.Item(23, 1).Value  = 3
$ws.Cells.Item(23, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(23, 3).Value  = "Coquimbo"
$ws.Cells.Item(23, 4).Value  = 44497
$ws.Cells.Item(23, 5).Value  = 5
$ws.Cells.Item(23, 6).Value  = 100112026
$ws.Cells.Item(23, 7).Value  = "Haba"
$ws.Cells.Item(23, 8).Value  = "Sin especificar"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 60
$ws.Cells.Item(23, 11).Value = 8000
$ws.Cells.Item(23, 12).Value = 8000
$ws.Cells.Item(23, 13).Value = 8000
$ws.Cells.Item(23, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(23, 16).Value = 320
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
